$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column R (18) into column S (19) for rows 3-34
$ws.Range("R3:R34").Copy()
$ws.Range("S3:S34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(4, 19).Value = 2022
$ws.Cells.Item(5, 19).Value = 135
$ws.Cells.Item(6, 19).Value = 99
$ws.Cells.Item(7, 19).Value = 36
$ws.Cells.Item(8, 19).Value = 97
$ws.Cells.Item(9, 19).Value = 80
$ws.Cells.Item(10, 19).Value = 17
$ws.Cells.Item(11, 19).Value = 17
$ws.Cells.Item(12, 19).Value = 11
$ws.Cells.Item(13, 19).Value = 6
$ws.Cells.Item(14, 19).Value = 5
$ws.Cells.Item(15, 19).Value = 3
$ws.Cells.Item(16, 19).Value = 2
$ws.Cells.Item(17, 19).Value = "-"
$ws.Cells.Item(18, 19).Value = "-"
$ws.Cells.Item(19, 19).Value = "-"
$ws.Cells.Item(20, 19).Value = 6
$ws.Cells.Item(21, 19).Value = 1
$ws.Cells.Item(22, 19).Value = 5
$ws.Cells.Item(23, 19).Value = "-"
$ws.Cells.Item(24, 19).Value = "-"
$ws.Cells.Item(25, 19).Value = "-"
$ws.Cells.Item(26, 19).Value = 10
$ws.Cells.Item(27, 19).Value = 4
$ws.Cells.Item(28, 19).Value = 6
$ws.Cells.Item(29, 19).Value = "-"
$ws.Cells.Item(30, 19).Value = "-"
$ws.Cells.Item(31, 19).Value = "-"
$ws.Cells.Item(32, 19).Value = "-"
$ws.Cells.Item(33, 19).Value = "-"
$ws.Cells.Item(34, 19).Value = "-"

$ws.Range("T24").Select()
